$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.339.51"
$ws.Range("E2").Value = "  +3.45%  "
$ws.Range("D3").Value = "2.271.19"
$ws.Range("E3").Value = "  +2.79%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "323.07"
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("D6").Value = "105.15"
$ws.Range("E6").Value = "  +6.32%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.571"
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.23%  "
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "7.88"
$ws.Range("E12").Value = "  +3.07%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("D15").Value = "2.618.64"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").Value = "14.58"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").Value = "2.275.80"
$ws.Range("E17").Value = "  +3.11%  "
$ws.Range("D18").Value = "44.252.83"
$ws.Range("E18").Value = "  +3.58%  "
$ws.Range("D19").Value = "13.87"
$ws.Range("E19").Value = "  -3.10%  "
$ws.Range("E20").Value = "  +4.65%  "
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").Value = "66.39"
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("D24").Value = "240.57"
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("E25").Value = "  +4.76%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "10.29"
$ws.Range("E27").Value = "  +3.40%  "
$ws.Range("D28").Value = "38.44"
$ws.Range("E28").Value = "  +12.53%  "
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("E30").Value = "  +3.41%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "163.06"
$ws.Range("E31").Value = "  +5.99%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "20.68"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("E35").Value = "  +9.05%  "
$ws.Range("E36").Value = "  +5.12%  "
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("E39").Value = "  +4.64%  "
$ws.Range("D40").Value = "4.41"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +1.87%  "
$ws.Range("D42").Value = "15.59"
$ws.Range("E42").Value = "  +28.24%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "1.782.18"
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("D46").Value = "86.25"
$ws.Range("E46").Value = "  -2.23%  "
$ws.Range("D47").Value = "5.42"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").Value = "60.68"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "1.73"
$ws.Range("E49").Value = "  +10.94%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "75.63"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").Value = "104.19"
